# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-30, ~17:55-17:56) to the
# PIR, Humidity, Proximity, mmWave(BR), mmWave(HR) and mmWave(InBed) sheets.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows,
        [bool]$TextifyColE
    )

    $ws = $wb.Worksheets.Item($SheetName)

    for ($i = 0; $i -lt $Rows.Length; $i++) {
        $r   = $StartRow + $i
        $row = $Rows[$i]

        # Column A holds dates formatted like "2026-01-30" - force Text so
        # Excel doesn't auto-convert the literal into a date serial.
        $ws.Cells.Item($r, 1).NumberFormat = "@"
        $ws.Cells.Item($r, 1).Value = $row[0]

        # Column B/C are plain HH:MM:SS / HH:MM strings - safe as-is.
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]

        # Column D is always a plain location/device label.
        $ws.Cells.Item($r, 4).Value = $row[3]

        # Column E is numeric on some sheets (mmWave BR/HR) and text
        # (including percentages like "86.7%") on others.
        if ($TextifyColE) {
            $ws.Cells.Item($r, 5).NumberFormat = "@"
        }
        $ws.Cells.Item($r, 5).Value = $row[4]

        # Column F is always a plain status/label string.
        $ws.Cells.Item($r, 6).Value = $row[5]
    }
}

# ---------------------------------------------------------------------------
# PIR sheet: 13 new "No Motion" / "Inactive" Bathroom readings (rows 389-401)
# ---------------------------------------------------------------------------
$pirRows = @(
    @("2026-01-30","17:55:12","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:14","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:17","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:22","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:27","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:32","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:37","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:42","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:47","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:53","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:55:58","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:03","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:08","17:00","Bathroom","No Motion","Inactive")
)
Append-Rows "PIR" 389 $pirRows $false

# ---------------------------------------------------------------------------
# Humidity sheet: 7 new Bathroom % readings (rows 264-270)
# ---------------------------------------------------------------------------
$humidityRows = @(
    @("2026-01-30","17:55:13","17:00","Bathroom","86.7%","Active"),
    @("2026-01-30","17:55:18","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:55:28","17:00","Bathroom","86.7%","Active"),
    @("2026-01-30","17:55:44","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:55:48","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:56:03","17:00","Bathroom","85.4%","Active"),
    @("2026-01-30","17:56:08","17:00","Bathroom","86.9%","Active")
)
Append-Rows "Humidity" 264 $humidityRows $true

# ---------------------------------------------------------------------------
# Proximity sheet: 9 new Bathroom Door ENTER/EXIT events (rows 65-73)
# ---------------------------------------------------------------------------
$proximityRows = @(
    @("2026-01-30","17:55:11","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:55:14","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:55:22","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:55:31","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:55:35","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:55:41","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:55:43","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:55:55","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:56:05","17:00","Bathroom Door","EXIT","User EXITED Bathroom")
)
Append-Rows "Proximity" 65 $proximityRows $false

# ---------------------------------------------------------------------------
# mmWave(BR): single new Bedroom reading, numeric Value = 0 (row 2)
# ---------------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")
$wsBR.Cells.Item(2, 1).NumberFormat = "@"
$wsBR.Cells.Item(2, 1).Value = "2026-01-30"
$wsBR.Cells.Item(2, 2).Value = "17:55:21"
$wsBR.Cells.Item(2, 3).Value = "17:00"
$wsBR.Cells.Item(2, 4).Value = "Bedroom"
$wsBR.Cells.Item(2, 5).Value = 0
$wsBR.Cells.Item(2, 6).Value = "Empty"

# ---------------------------------------------------------------------------
# mmWave(HR): single new Bedroom reading, numeric Value = 0 (row 2)
# ---------------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")
$wsHR.Cells.Item(2, 1).NumberFormat = "@"
$wsHR.Cells.Item(2, 1).Value = "2026-01-30"
$wsHR.Cells.Item(2, 2).Value = "17:55:20"
$wsHR.Cells.Item(2, 3).Value = "17:00"
$wsHR.Cells.Item(2, 4).Value = "Bedroom"
$wsHR.Cells.Item(2, 5).Value = 0
$wsHR.Cells.Item(2, 6).Value = "Empty"

# ---------------------------------------------------------------------------
# mmWave(InBed): single new Bedroom reading, text Value = "Out of Bed" (row 2)
# ---------------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")
$wsInBed.Cells.Item(2, 1).NumberFormat = "@"
$wsInBed.Cells.Item(2, 1).Value = "2026-01-30"
$wsInBed.Cells.Item(2, 2).Value = "17:55:19"
$wsInBed.Cells.Item(2, 3).Value = "17:00"
$wsInBed.Cells.Item(2, 4).Value = "Bedroom"
$wsInBed.Cells.Item(2, 5).Value = "Out of Bed"
$wsInBed.Cells.Item(2, 6).Value = "Empty"
